$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Split the opening paragraph into a greeting line ("Dear Captain,")
#    and a separate body paragraph with a larger first-line indent.
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("     I had learn a lot in this course.", $true, $false, $false, $false, $false, $true, 1, $false, "     Dear Captain,^pI had learn a lot in this course.", 2) | Out-Null

# Locate the newly created "I had learn a lot..." paragraph and give it
# the larger indent used for the greeting's body paragraph.
$find2 = $d.Content.Find
$find2.Execute("I had learn a lot in this course. This class let me remember", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bodyPara = $find2.Parent.Paragraphs.Item(1)
$bodyPara.Format.CharacterUnitFirstLineIndent = 350
$bodyPara.Format.FirstLineIndent = 36.75

# ------------------------------------------------------------------
# 2. Closing paragraph edits.
# ------------------------------------------------------------------

# 2a. "teacher" -> "you and Miss Shannon" right before "can give us feedback"
$d.Content.Find.Execute("next time teacher can give", $true, $false, $false, $false, $false, $true, 1, $false, "next time you and Miss Shannon can give", 2) | Out-Null

# 2b. "draft2" -> "draft 2" in "doesn't drag the draft2 to"
$d.Content.Find.Execute("doesn't drag the draft2 to", $true, $false, $false, $false, $false, $true, 1, $false, "doesn't drag the draft 2 to", 2) | Out-Null

# 2c. Replace the "I hope teacher..." sentence with the new thanks/hope sentence.
$d.Content.Find.Execute("of the course. I hope teacher could give us the draft2 on time so we won`u2019t have too much works to do at the end of the class. I think I think next time", $true, $false, $false, $false, $false, $true, 1, $false, "of the course. I also hope we won`u2019t have many technique problem, so we can have our class immediately. I think I think next time", 2) | Out-Null

# 2d. Append the closing thank-you sentence and a trailing tab.
$d.Content.Find.Execute("I would still try my best to work on homework and classes.", $true, $false, $false, $false, $false, $true, 1, $false, "I would still try my best to work on homework and classes. Thank you Captain and Miss Shannon!^t", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Remove the _GoBack bookmark from its old position (mid-paragraph);
#    it will be re-added at the end of the document below.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 4. Append the signature block as new paragraphs after the closing
#    paragraph, each with the same indentation style used elsewhere.
# ------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.Move(1, 1) | Out-Null
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "                                                   Your Student,`r"

$endRange2 = $d.Content
$endRange2.Collapse(0)
$p2 = $d.Paragraphs.Last
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "                                                      Michelle Wu`r"

$p4 = $d.Paragraphs.Last
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$p5.Range.Text = "                                                                      2020/8/16"

foreach ($idx in 9..11) {
    $p = $d.Paragraphs.Item($idx)
    $p.Format.CharacterUnitFirstLineIndent = 250
    $p.Format.FirstLineIndent = 26.25
}

$lastRange = $d.Paragraphs.Last.Range
$lastRange.Collapse(0)
$lastRange.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $lastRange) | Out-Null

Write-Output "done"
